$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (index row) updates
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON) updates
$ws.Range("B2").Value = 389.61793576387959
$ws.Range("C2").Value = 459.71363799849394
$ws.Range("D2").Value = 385.30257189742673
$ws.Range("E2").Value = 464.06631234314534

# Row 3 (STR) updates
$ws.Range("B3").Value = 388.3339396377682
$ws.Range("C3").Value = 476.53000135572012
$ws.Range("D3").Value = 395.4634356969558
$ws.Range("E3").Value = 468.18968998464879

# Selection change
$ws.Range("B1:E3").Select() | Out-Null
